$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly fruit/vegetable data reshuffle: rows 2-28 (excluding 16-17, which are
# unchanged) had their Fecha / Volumen / Precio* / Origen / Calidad values moved
# around to reflect the corrected weekly ordering.

$ws.Range("D2").Value = 44161
$ws.Range("J2").Value = 1600

$ws.Range("D3").Value = 44161
$ws.Range("J3").Value = 1850

$ws.Range("D4").Value = 44174
$ws.Range("J4").Value = 2800
$ws.Range("K4").Value = 1200
$ws.Range("L4").Value = 1250
$ws.Range("M4").Value = 1221
$ws.Range("P4").Value = 1221

$ws.Range("D5").Value = 44174
$ws.Range("J5").Value = 1300

$ws.Range("D6").Value = 44160
$ws.Range("J6").Value = 750

$ws.Range("D7").Value = 44160
$ws.Range("J7").Value = 850

$ws.Range("D8").Value = 44159
$ws.Range("J8").Value = 1100
$ws.Range("K8").Value = 1300
$ws.Range("L8").Value = 1300
$ws.Range("M8").Value = 1300
$ws.Range("P8").Value = 1300

$ws.Range("D9").Value = 44159
$ws.Range("J9").Value = 800

$ws.Range("D10").Value = 44169
$ws.Range("J10").Value = 950

$ws.Range("D11").Value = 44169
$ws.Range("J11").Value = 800
$ws.Range("K11").Value = 1000
$ws.Range("L11").Value = 1000
$ws.Range("M11").Value = 1000
$ws.Range("P11").Value = 1000

$ws.Range("D12").Value = 44181
$ws.Range("J12").Value = 1000
$ws.Range("K12").Value = 1300
$ws.Range("M12").Value = 1300
$ws.Range("P12").Value = 1300

$ws.Range("D13").Value = 44181
$ws.Range("J13").Value = 900
$ws.Range("K13").Value = 900
$ws.Range("L13").Value = 900
$ws.Range("M13").Value = 900
$ws.Range("P13").Value = 900

$ws.Range("D14").Value = 44175
$ws.Range("J14").Value = 1500

$ws.Range("D15").Value = 44175
$ws.Range("J15").Value = 1450

$ws.Range("D18").Value = 44165
$ws.Range("J18").Value = 720
$ws.Range("K18").Value = 1200
$ws.Range("L18").Value = 1200
$ws.Range("M18").Value = 1200
$ws.Range("P18").Value = 1200

$ws.Range("D19").Value = 44165
$ws.Range("J19").Value = 750

$ws.Range("D20").Value = 44167
$ws.Range("J20").Value = 1430
$ws.Range("K20").Value = 1200
$ws.Range("M20").Value = 1248
$ws.Range("P20").Value = 1248

$ws.Range("D21").Value = 44167
$ws.Range("J21").Value = 350

$ws.Range("D22").Value = 44176
$ws.Range("J22").Value = 2500
$ws.Range("L22").Value = 1300
$ws.Range("M22").Value = 1256
$ws.Range("P22").Value = 1256

$ws.Range("D23").Value = 44176
$ws.Range("J23").Value = 1500

$ws.Range("D24").Value = 44179
$ws.Range("J24").Value = 980
$ws.Range("L24").Value = 1200
$ws.Range("M24").Value = 1200
$ws.Range("O24").Value = "Región Metropolitana"
$ws.Range("P24").Value = 1200

$ws.Range("D25").Value = 44162
$ws.Range("I25").Value = "Primera"
$ws.Range("J25").Value = 1200
$ws.Range("K25").Value = 1300
$ws.Range("L25").Value = 1300
$ws.Range("M25").Value = 1300
$ws.Range("P25").Value = 1300

$ws.Range("D26").Value = 44162
$ws.Range("I26").Value = "Segunda"
$ws.Range("J26").Value = 800
$ws.Range("K26").Value = 1000
$ws.Range("L26").Value = 1000
$ws.Range("M26").Value = 1000
$ws.Range("O26").Value = "Provincia de Quillota"
$ws.Range("P26").Value = 1000

$ws.Range("D27").Value = 44168
$ws.Range("J27").Value = 1200

$ws.Range("D28").Value = 44168
